$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.07"
$ws.Range("E2").Value = "'-0.43%"
$ws.Range("E3").Value = "'3.40%"
$ws.Range("D4").Value = "'5.149"
$ws.Range("E4").Value = "'0.92%"
$ws.Range("E5").Value = "'0.33%"
$ws.Range("D6").Value = "'6.466"
$ws.Range("E6").Value = "'-0.11%"
$ws.Range("D7").Value = "'0.8198"
$ws.Range("E7").Value = "'0.17%"
$ws.Range("D8").Value = "'0.8287"
$ws.Range("E8").Value = "'-1.43%"
$ws.Range("E9").Value = "'-0.55%"
$ws.Range("D11").Value = "'0.02886"
$ws.Range("E11").Value = "'0.69%"
$ws.Range("E12").Value = "'0.06%"
$ws.Range("D13").Value = "'0.001506"
$ws.Range("E13").Value = "'-0.68%"
$ws.Range("D14").Value = "'0.0005966"
$ws.Range("E14").Value = "'0.10%"
$ws.Range("D15").Value = "'0.006228"
$ws.Range("E15").Value = "'0.36%"
$ws.Range("D16").Value = "'3.654"
$ws.Range("E16").Value = "'3.72%"
$ws.Range("D17").Value = "'3.034"
$ws.Range("E17").Value = "'0.62%"
$ws.Range("E18").Value = "'-12.66%"
$ws.Range("E19").Value = "'-2.12%"
$ws.Range("D22").Value = "'3.741"
$ws.Range("E22").Value = "'-0.18%"
$ws.Range("E23").Value = "'-1.78%"
$ws.Range("E24").Value = "'-2.48%"
$ws.Range("D25").Value = "'0.001223"
$ws.Range("E25").Value = "'-1.87%"
$ws.Range("E26").Value = "'-2.75%"
$ws.Range("D27").Value = "'0.00009597"
$ws.Range("E27").Value = "'-1.07%"
$ws.Range("D28").Value = "'0.0001398"
$ws.Range("E28").Value = "'0.60%"
$ws.Range("D40").Value = "'0.03645"
$ws.Range("E40").Value = "'-0.41%"
$ws.Range("D41").Value = "'0.1368"
$ws.Range("E41").Value = "'30.08%"
$ws.Range("D42").Value = "'0.006158"
$ws.Range("E42").Value = "'-0.04%"
$ws.Range("E43").Value = "'4.02%"
$ws.Range("D44").Value = "'0.009015"
$ws.Range("E44").Value = "'7.12%"
$ws.Range("D45").Value = "'0.00005346"
$ws.Range("E45").Value = "'0.47%"
$ws.Range("E46").Value = "'-0.08%"
$ws.Range("E47").Value = "'8.18%"
$ws.Range("D48").Value = "'0.002336"
$ws.Range("E48").Value = "'10.04%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.08%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'-0.08%"
